# Add 2022-Q1 data:
#  1) insert a new "2022-Q1" worksheet (same layout as the other quarterly
#     sheets) right before the "总计" (totals) sheet
#  2) update the "总计" sheet with a new top row for 2022-Q1 and correct the
#     2021-Q3 holding value

$wb = $excel.ActiveWorkbook

# --- Step 1: insert the "2022-Q1" worksheet -------------------------------
$sheetTotal = $wb.Worksheets.Item(6)          # "总计" sheet (currently last)
$refSheet   = $wb.Worksheets.Item(5)          # "2021-Q4" sheet, used as a style template

$newSheet = $wb.Worksheets.Add($sheetTotal)   # inserts the new sheet right before "总计"
$newSheet.Name = "2022-Q1"

# copy the header-row and index-cell formatting from the reference sheet
$refSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats
$refSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)     # xlPasteFormats

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("B2:G2").NumberFormat = "@"   # keep numeric-looking values as text
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "007280"
$newSheet.Range("C2").Value = "上投摩根日本精选股票（QDII）"
$newSheet.Range("D2").Value = "1.35"
$newSheet.Range("E2").Value = "88.71"
$newSheet.Range("F2").Value = "2.53"
$newSheet.Range("G2").Value = "0.0342"
$newSheet.Range("H2").Value = 7
$newSheet.Range("B2:G2").Style = "Normal"

# --- Step 2: update the "总计" summary sheet -------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()             # shift existing data rows down by one

# match the formatting of the (now shifted) index column
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.03
$totalSheet.Range("B2:D2").Style = "Normal"

# renumber the index column (A) for the rows that got shifted down
for ($r = 3; $r -le 7; $r++) {
    $totalSheet.Range("A$r").Value = $r - 2
}

# correct the 2021-Q3 holding value (now on row 4): 0.04 -> 0.03
$totalSheet.Range("D4").Value = 0.03
